# Update the date line and all twenty-five multiplication equations
# in the table to the new values per the commit diff.

$d = $word.ActiveDocument

$replacements = @(
    @("2024-09-30 Monday", "2024-10-01 Tuesday"),
    @("787×7=5509", "109×4=436"),
    @("346×2=692", "553×3=1659"),
    @("467×4=1868", "453×9=4077"),
    @("780×7=5460", "123×4=492"),
    @("208×6=1248", "730×3=2190"),
    @("910×6=5460", "507×3=1521"),
    @("317×9=2853", "394×2=788"),
    @("547×7=3829", "895×5=4475"),
    @("652×2=1304", "899×6=5394"),
    @("406×5=2030", "680×3=2040"),
    @("491×8=3928", "633×4=2532"),
    @("985×4=3940", "599×6=3594"),
    @("356×8=2848", "784×6=4704"),
    @("435×4=1740", "347×3=1041"),
    @("447×9=4023", "232×6=1392"),
    @("506×8=4048", "229×4=916"),
    @("354×6=2124", "337×2=674"),
    @("977×8=7816", "592×9=5328"),
    @("411×6=2466", "816×5=4080"),
    @("658×8=5264", "721×9=6489"),
    @("869×6=5214", "533×4=2132"),
    @("730×7=5110", "563×2=1126"),
    @("506×7=3542", "498×8=3984"),
    @("497×8=3976", "148×8=1184"),
    @("841×3=2523", "130×2=260")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
